$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.ShowAllData()
